$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # DataFetchFlag
$ws2 = $wb.Worksheets.Item(2)  # DataFetchXL

# --- Sheet "DataFetchFlag": append new row 22 ---
$ws1.Range("A22").Value = "EXL_CorporateLensHomePage_Announcements"
$ws1.Range("B22").Value = "XL"

# Re-apply the list data validation over the extended range B2:B22
$ws1.Range("B2:B21").Validation.Delete()
$ws1.Range("B2:B22").Validation.Add(3, 1, 1, '"XL,DB"')

# --- Sheet "DataFetchXL": append new row 22 ---
$ws2.Range("A22").Value = "EXL_CorporateLensHomePage_Announcements"
$ws2.Range("B22").Value = "\\src\\com\\proj\\suitecorporateLens\\testdata\\CorporateLensTestData-Announcements.xlsx "
$ws2.Range("C22").Value = "Announcements"

# Hyperlink the new path cell, same as the existing rows above it
$ws2.Hyperlinks.Add($ws2.Range("B22"), "file:///\\src\com\proj\suitecorporateLens\testdata\CorporateLensTestData-Announcements.xlsx ")
$ws2.Range("B22").Style = "Hyperlink"

# Update the visible selection to match the newly entered cell on each sheet
$ws2.Range("C22").Select()

# Re-select sheet1's new cell last so it remains the active/visible tab
# (matches the workbook's original active-sheet state)
$ws1.Range("B22").Select()
